$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 63003.5585
$ws.Range("C2").Value = 5373.441499999999
$ws.Range("D2").Value = 45000
$ws.Range("E2").Value = 12630.117
$ws.Range("B3").Value = 60503.643
$ws.Range("C3").Value = 5250.357
$ws.Range("D3").Value = 22500
$ws.Range("E3").Value = 32753.28599999999
$ws.Range("B4").Value = 58280.113
$ws.Range("C4").Value = 5171.886999999999
$ws.Range("D4").Value = 22500
$ws.Range("E4").Value = 30608.226
$ws.Range("B5").Value = 56670.0195
$ws.Range("C5").Value = 5107.980499999999
$ws.Range("D5").Value = 22500
$ws.Range("E5").Value = 29062.039
$ws.Range("B6").Value = 58596.607
$ws.Range("C6").Value = 5104.393
$ws.Range("D6").Value = 22500
$ws.Range("E6").Value = 30992.21400000001
$ws.Range("B7").Value = 61690.523
$ws.Range("C7").Value = 5230.476999999999
$ws.Range("D7").Value = 22500
$ws.Range("E7").Value = 33960.046
$ws.Range("B8").Value = 59656.1205
$ws.Range("C8").Value = 5902.879499999999
$ws.Range("D8").Value = 22500
$ws.Range("E8").Value = 31253.24099999999
$ws.Range("B9").Value = 48688.998
$ws.Range("C9").Value = 7037.001999999999
$ws.Range("D9").Value = 45000
$ws.Range("E9").Value = -3348.004000000001
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 8592.4405
$ws.Range("D10").Value = 65000
$ws.Range("E10").ClearContents()
$ws.Range("B11").Value = 90861.24603960395
$ws.Range("C11").Value = 14095.0215
$ws.Range("D11").Value = 65000
$ws.Range("E11").Value = 11766.22453960395
$ws.Range("B12").Value = 89726.60198019801
$ws.Range("C12").Value = 15828.8445
$ws.Range("D12").Value = 65000
$ws.Range("E12").Value = 8897.757480198023
$ws.Range("B13").Value = 90887.15396039604
$ws.Range("C13").Value = 15449.287
$ws.Range("D13").Value = 65000
$ws.Range("E13").Value = 10437.86696039604
$ws.Range("B14").Value = 95761.86584158416
$ws.Range("C14").Value = 15803.2595
$ws.Range("D14").Value = 65000
$ws.Range("E14").Value = 14958.60634158416
$ws.Range("B15").Value = 74053.73267326732
$ws.Range("C15").Value = 15862.028
$ws.Range("D15").Value = 65000
$ws.Range("E15").Value = -6808.295326732674
$ws.Range("B16").Value = 95771.50148514852
$ws.Range("C16").Value = 16032.3205
$ws.Range("D16").Value = 65000
$ws.Range("E16").Value = 14739.18098514852
$ws.Range("B17").Value = 84547.68217821782
$ws.Range("C17").Value = 16559.4275
$ws.Range("D17").Value = 65000
$ws.Range("E17").Value = 2988.254678217811
$ws.Range("B18").Value = 59953.5995049505
$ws.Range("C18").Value = 16683.57249999999
$ws.Range("D18").Value = 65000
$ws.Range("E18").Value = -21729.97299504949
$ws.Range("B19").Value = 72066.9185
$ws.Range("C19").Value = 16252.82049999999
$ws.Range("D19").Value = 65000
$ws.Range("E19").Value = -9185.901999999987
$ws.Range("B20").Value = 74173.3315
$ws.Range("C20").Value = 15546.958
$ws.Range("D20").Value = 65000
$ws.Range("E20").Value = -6373.626499999998
$ws.Range("B21").Value = 72400.7075
$ws.Range("C21").Value = 13748.4165
$ws.Range("D21").Value = 65000
$ws.Range("E21").Value = -6347.708999999995
$ws.Range("B22").Value = 70217.473
$ws.Range("C22").Value = 12030.7215
$ws.Range("D22").Value = 65000
$ws.Range("E22").Value = -6813.248500000002
$ws.Range("B23").Value = 21309.33
$ws.Range("C23").Value = 9319.645999999999
$ws.Range("D23").Value = 65000
$ws.Range("E23").Value = -53010.316
$ws.Range("B24").Value = 69506.016
$ws.Range("C24").Value = 6308.333499999995
$ws.Range("D24").Value = 65000
$ws.Range("E24").Value = -1802.31749999999
$ws.Range("B25").Value = 44605.08100000001
$ws.Range("C25").Value = 5687.261999999999
$ws.Range("D25").Value = 57000
$ws.Range("E25").Value = -18082.18099999998
